$wb = $excel.ActiveWorkbook

# --- 1) Update status text "Ready for handoff" -> "In Translation" ---
# Overview sheet: columns E (zh-cn) and F (de-de) in row 2
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

# zh-cn sheet: Status column C, row 2
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

# de-de sheet: Status column C, row 2
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# --- 2) Narrow the date/status-related columns ---
# Overview sheet: columns E and F width 17.22 -> 13.41 (characters)
$wsOverview.Columns.Item(5).ColumnWidth = 12.43
$wsOverview.Columns.Item(6).ColumnWidth = 12.43

# zh-cn / de-de sheets: column C width 17.22 -> 13.41 (characters)
$wsZhCn.Columns.Item(3).ColumnWidth = 12.43
$wsDeDe.Columns.Item(3).ColumnWidth = 12.43
